$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) both contain the same
# event rows in columns F2 and F3 ("想去人数" - number of people who
# want to go) that need to be refreshed with newly scraped counts.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 6602
    $ws.Range("F3").Value = 41
}
